$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated nearest-hospital results (Policy Maker Model)
$ws.Range("D3").Value  = 49.26812930661734
$ws.Range("E3").Value  = 6.929470619593767
$ws.Range("F3").Value  = 2.377

$ws.Range("D7").Value  = 49.495769
$ws.Range("E7").Value  = 6.596771
$ws.Range("F7").Value  = 10.059

$ws.Range("D10").Value = 49.457186
$ws.Range("E10").Value = 6.631578
$ws.Range("F10").Value = 10.122

$ws.Range("D11").Value = 49.537868
$ws.Range("E11").Value = 6.88934
$ws.Range("F11").Value = 4.441

$ws.Range("D12").Value = 49.33767769818895
$ws.Range("E12").Value = 7.005279273288449
$ws.Range("F12").Value = 3.439

$ws.Range("D13").Value = 49.33767769818895
$ws.Range("E13").Value = 7.005279273288449
$ws.Range("F13").Value = 6.382

$ws.Range("D14").Value = 49.40973615864516
$ws.Range("E14").Value = 7.172516962777409
$ws.Range("F14").Value = 4.915

$ws.Range("D15").Value = 49.40973615864516
$ws.Range("E15").Value = 7.172516962777409
$ws.Range("F15").Value = 2.272

$ws.Range("D16").Value = 49.39718981222195
$ws.Range("E16").Value = 7.213364346793378
$ws.Range("F16").Value = 2.336

$ws.Range("D20").Value = 49.35666769333056
$ws.Range("E20").Value = 6.822620825380222
$ws.Range("F20").Value = 1.204

$ws.Range("D21").Value = 49.44208949482798
$ws.Range("E21").Value = 6.904224529034402
$ws.Range("F21").Value = 1.52

$ws.Range("D26").Value = 49.3762045599062
$ws.Range("E26").Value = 7.280034712412421
$ws.Range("F26").Value = 5.111

$ws.Range("D28").Value = 49.537868
$ws.Range("E28").Value = 6.88934
$ws.Range("F28").Value = 10.326

$ws.Range("D29").Value = 49.453979
$ws.Range("E29").Value = 7.178492
$ws.Range("F29").Value = 9.693

$ws.Range("D30").Value = 49.454659
$ws.Range("E30").Value = 7.186793
$ws.Range("F30").Value = 7.876

$ws.Range("D31").Value = 49.454659
$ws.Range("E31").Value = 7.186793
$ws.Range("F31").Value = 12.617
